# registration.xlsx — implement user permissions for the "properties" sheet:
#   - change the existing Table/security/filterTypeOnCreation row's value
#     from "MODIFY" to "HIDDEN"
#   - add a new Table/security/locked row of type boolean, value true
#   - leave a pre-formatted (quote-prefixed) blank cell below it (E4)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("properties")

# Existing filterTypeOnCreation row: MODIFY -> HIDDEN
$ws.Range("E2").Value = "HIDDEN"

# New "locked" security row
$ws.Range("A3").Value = "Table"
$ws.Range("B3").Value = "security"
$ws.Range("C3").Value = "locked"
$ws.Range("D3").Value = "boolean"
$ws.Range("E3").Value = "'true"

# Leftover quote-prefixed (text) formatting on the next row's E cell, no value
$ws.Range("E4").Value = "'true"
$ws.Range("E4").ClearContents()

$ws.Activate()
$ws.Range("J29").Select()
